$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C16").Value = "1007738481"
$ws.Range("D16").Value = "LUIS ANTONIO DE AVILA HERNANDEZ"
$ws.Range("E16").Value = "2212"
$ws.Range("F16").Value = 6667
$ws.Range("G16").Value = 1000000

$ws.Range("C17").Value = "1082856017"
$ws.Range("D17").Value = "RUBEN DARIO CUCUNUBA SALINAS"
$ws.Range("E17").Value = "2301"
$ws.Range("F17").Value = 41760
$ws.Range("G17").Value = 1160000

$ws.Range("C18").Value = "1082856017"
$ws.Range("D18").Value = "RUBEN DARIO CUCUNUBA SALINAS"
$ws.Range("E18").Value = "2301"
$ws.Range("F18").Value = 40000
$ws.Range("G18").Value = 1000000

$ws.Range("C19").Value = "1047376325"
$ws.Range("D19").Value = "JUAN DAVID PAEZ CORTECERO"
$ws.Range("E19").Value = "2301"
$ws.Range("F19").Value = 27840
$ws.Range("G19").Value = 1160000

$ws.Range("C20").Value = "1082856017"
$ws.Range("D20").Value = "RUBEN DARIO CUCUNUBA SALINAS"
$ws.Range("E20").Value = "2302"
$ws.Range("F20").Value = 46400
$ws.Range("G20").Value = 1160000

$ws.Range("C21").Value = "1047376325"
$ws.Range("D21").Value = "JUAN DAVID PAEZ CORTECERO"
$ws.Range("E21").Value = "2302"
$ws.Range("F21").Value = 46400
$ws.Range("G21").Value = 1160000

$ws.Range("C22").Value = "1082856017"
$ws.Range("D22").Value = "RUBEN DARIO CUCUNUBA SALINAS"
$ws.Range("E22").Value = "2303"
$ws.Range("F22").Value = 46400
$ws.Range("G22").Value = 1160000

$ws.Range("C23").Value = "1050944608"
$ws.Range("D23").Value = "ALBEIRO ENRIQUE PEREZ CAUSIL"
$ws.Range("E23").Value = "2303"
$ws.Range("F23").Value = 34027
$ws.Range("G23").Value = 1160000

$ws.Range("C24").Value = "1047376325"
$ws.Range("D24").Value = "JUAN DAVID PAEZ CORTECERO"
$ws.Range("E24").Value = "2303"
$ws.Range("F24").Value = 46400
$ws.Range("G24").Value = 1160000

$ws.Range("C25").Value = "1050967361"
$ws.Range("D25").Value = "JUAN CAMILO HERNANDEZ VIGGIANI"
$ws.Range("E25").Value = "2304"
$ws.Range("F25").Value = 46400
$ws.Range("G25").Value = 1160000

$ws.Range("C26").Value = "1066734978"
$ws.Range("D26").Value = "JHON DEIVY GARCIA SANCHEZ"
$ws.Range("E26").Value = "2304"
$ws.Range("F26").Value = 46400
$ws.Range("G26").Value = 1160000

$ws.Range("C27").Value = "1050967361"
$ws.Range("D27").Value = "JUAN CAMILO HERNANDEZ VIGGIANI"
$ws.Range("E27").Value = "2305"
$ws.Range("F27").Value = 37120
$ws.Range("G27").Value = 1160000

$ws.Range("C28").Value = "1066734978"
$ws.Range("D28").Value = "JHON DEIVY GARCIA SANCHEZ"
$ws.Range("E28").Value = "2305"
$ws.Range("F28").Value = 37120
$ws.Range("G28").Value = 1160000
